$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells so numeric-looking strings
# (e.g. "1.00", "0.999") are preserved exactly instead of being
# parsed into numbers. (Set one cell at a time: a combined
# multi-area Range(...).NumberFormat assignment only took effect
# on the first area when tested against this runtime.)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '91.328.64'
$ws.Range("E2").Value = '  +2.55%  '
$ws.Range("D3").Value = '3.118.92'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '218.87'
$ws.Range("E5").Value = '  +2.77%  '
$ws.Range("D6").Value = '623.10'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '1.02'
$ws.Range("E7").Value = '  +25.83%  '
$ws.Range("D8").Value = '0.375'
$ws.Range("E8").Value = '  -0.41%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").Value = '3.114.51'
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("D11").Value = '0.736'
$ws.Range("E11").Value = '  +23.29%  '
$ws.Range("D12").Value = '0.193'
$ws.Range("E12").Value = '  +6.29%  '
$ws.Range("D13").Value = '0.0000253'
$ws.Range("E13").Value = '  +4.45%  '
$ws.Range("D14").Value = '34.68'
$ws.Range("E14").Value = '  +7.21%  '
$ws.Range("D15").Value = '5.49'
$ws.Range("E15").Value = '  +3.67%  '
$ws.Range("D16").Value = '91.115.26'
$ws.Range("E16").Value = '  +2.91%  '
$ws.Range("D17").Value = '3.690.37'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("B18").Value = 'SuiNetwork'
$ws.Range("C18").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D18").Value = '3.88'
$ws.Range("E18").Value = '  +14.60%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.101.65'
$ws.Range("E19").Value = '  +0.69%  '
$ws.Range("D20").Value = '0.0000218'
$ws.Range("E20").Value = '  +2.96%  '
$ws.Range("D21").Value = '14.12'
$ws.Range("E21").Value = '  +4.79%  '
$ws.Range("D22").Value = '441.92'
$ws.Range("E22").Value = '  +4.19%  '
$ws.Range("D23").Value = '8.89'
$ws.Range("E23").Value = '  +7.24%  '
$ws.Range("E24").Value = '  +5.15%  '
$ws.Range("D25").Value = '6.20'
$ws.Range("E25").Value = '  +9.78%  '
$ws.Range("D26").Value = '88.83'
$ws.Range("E26").Value = '  +7.54%  '
$ws.Range("D27").Value = '12.33'
$ws.Range("E27").Value = '  +3.11%  '
$ws.Range("D28").Value = '3.281.31'
$ws.Range("E28").Value = '  +1.53%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '0.166'
$ws.Range("E30").Value = '  -2.14%  '
$ws.Range("D31").Value = '9.20'
$ws.Range("E31").Value = '  +13.06%  '
$ws.Range("D32").Value = '525.73'
$ws.Range("E32").Value = '  +2.65%  '
$ws.Range("D33").Value = '0.897'
$ws.Range("E33").Value = '  -16.66%  '
$ws.Range("D34").Value = '3.75'
$ws.Range("E34").Value = '  +1.25%  '
$ws.Range("E35").Value = '  +13.39%  '
$ws.Range("D36").Value = '7.09'
$ws.Range("E36").Value = '  +4.22%  '
$ws.Range("D37").Value = '23.94'
$ws.Range("E37").Value = '  +7.15%  '
$ws.Range("E38").Value = '  +3.41%  '
$ws.Range("D39").Value = '1.87'
$ws.Range("E39").Value = '  +3.30%  '
$ws.Range("D40").Value = '0.0871'
$ws.Range("E40").Value = '  +25.01%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").Value = '0.156'
$ws.Range("E43").Value = '  +17.76%  '
$ws.Range("D44").Value = '0.398'
$ws.Range("E44").Value = '  +9.09%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = '1.94'
$ws.Range("E45").Value = '  +5.86%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = '148.56'
$ws.Range("E47").Value = '  +1.53%  '
$ws.Range("D48").Value = '43.99'
$ws.Range("E48").Value = '  +1.18%  '
$ws.Range("E49").Value = '  +6.50%  '
$ws.Range("D50").Value = '169.00'
$ws.Range("E50").Value = '  +3.50%  '
$ws.Range("D51").Value = '4.25'
$ws.Range("E51").Value = '  +7.44%  '
